$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "Baza podataka" -> "Analiza sadržaja" in cell C3
$ws.Range("C3").Value = "Analiza sadržaja"

# Reflect the resulting selection (user clicked/edited C3)
$ws.Range("C3").Select()
